$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column A for "Sprint" - shifts existing columns A..L to B..M
$ws.Range("A1").EntireColumn.Insert()

# Row 1 headers
$ws.Range("A1").Value = "Sprint"
$ws.Range("B1").Value = "Nombre del analista"
$ws.Range("C1").Value = "Habilitadores"
$ws.Range("D1").Value = "Historias de Usuario"
$ws.Range("E1").Value = "Bug"
$ws.Range("F1").Value = "Estado ""New"""
$ws.Range("G1").Value = "Estado ""Active"""
$ws.Range("H1").Value = "Estado ""Closed"""
$ws.Range("I1").Value = "Estado ""Impedimento"""
# J1, K1, L1 keep their previous text (shifted from I1, J1, K1)
$ws.Range("J1").Value = "Comprometido a X historias de usuario"
$ws.Range("K1").Value = "No puntuadas"
$ws.Range("L1").Value = "Numero de pull requests"
$ws.Range("M1").Value = "Numero de commits"

# Row 2 data
$ws.Range("A2").Value = "97"
$ws.Range("B2").Value = "Andres Felipe Blandon Palacio"
$ws.Range("K2").Value = "0"
$ws.Range("M2").Value = 0

# Row 3 data
$ws.Range("A3").Value = "97"
$ws.Range("B3").Value = "Juan David Londono Agudelo"
$ws.Range("K3").Value = "0"
$ws.Range("M3").Value = 0

# Row 4 data
$ws.Range("A4").Value = "97"
$ws.Range("B4").Value = "Alex Alberto Franco Cano"
$ws.Range("K4").Value = "0"
$ws.Range("M4").Value = 0

# Apply header style (s=1) to new M1 cell, and data style (s=2) to new M column cells
$ws.Range("M1").Style = $ws.Range("L1").Style
$ws.Range("M2:M4").Style = $ws.Range("L2:L4").Style

# Column width for new column M (col 13) to match others
$ws.Columns("M").ColumnWidth = $ws.Columns("L").ColumnWidth
